# Applies the "adding power plants scenario and hydrogen demand TS" edit
# to the "Coupling Parameters" sheet of the workbook.

$wb = $excel.ActiveWorkbook

# Target the "Coupling Parameters" sheet explicitly (it is also the
# workbook's active sheet, but look it up by name to be robust).
$ws = $wb.Worksheets.Item("Coupling Parameters")

# B3: End Year 2065 -> 2070
$ws.Range("B3").Value = 2070

# Row 20 ("Look Ahead" row): shrink the custom row height from 29 to 15
$ws.Rows.Item(20).RowHeight = 15

# B30: maximum_investment_capacity_per_year 1,000,000 -> 100,000,000
$ws.Range("B30").Value = 100000000

# B31: new boolean switch FALSE -> TRUE
$ws.Range("B31").Value = $true

# Reflect the updated view/selection state (scrolled down a bit further,
# active cell on the newly toggled switch).
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 14
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("B31").Select()
